$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Formula = '="38.043.84"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E2').Value = '  +2.80%  '

$c = $ws.Range('D3')
$c.Formula = '="2.056.63"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E3').Value = '  +2.07%  '

$ws.Range('E4').Value = '  +0.18%  '

$c = $ws.Range('D5')
$c.Formula = '="230.33"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.90%  '

$ws.Range('E6').Value = '  +1.42%  '

$c = $ws.Range('D7')
$c.Formula = '="58.30"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E7').Value = '  +6.27%  '

$ws.Range('E8').Value = '  +0.01%  '

$c = $ws.Range('D9')
$c.Formula = '="0.386"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E9').Value = '  +2.68%  '

$c = $ws.Range('D10')
$c.Formula = '="0.0808"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E10').Value = '  +3.41%  '

$ws.Range('E11').Value = '  +0.67%  '

$c = $ws.Range('D12')
$c.Formula = '="2.361.39"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E12').Value = '  +2.12%  '

$c = $ws.Range('D13')
$c.Formula = '="14.60"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E13').Value = '  +3.38%  '

$c = $ws.Range('D14')
$c.Formula = '="20.65"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E14').Value = '  +2.20%  '

$c = $ws.Range('D15')
$c.Formula = '="0.753"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E15').Value = '  +1.63%  '

$c = $ws.Range('D16')
$c.Formula = '="5.27"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E16').Value = '  +2.97%  '

$c = $ws.Range('D17')
$c.Formula = '="2.059.85"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.86%  '

$c = $ws.Range('D18')
$c.Formula = '="37.954.58"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E18').Value = '  +2.60%  '

$ws.Range('E19').Value = '  -0.73%  '

$c = $ws.Range('D20')
$c.Formula = '="69.89"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E20').Value = '  +1.57%  '

$c = $ws.Range('D21')
$c.Formula = '="0.0₃0830"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E21').Value = '  +1.78%  '

$c = $ws.Range('D22')
$c.Formula = '="224.71"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E22').Value = '  +0.83%  '

$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('E24').Value = '  +1.23%  '

$c = $ws.Range('D25')
$c.Formula = '="2.24"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E25').Value = '  +2.92%  '

$ws.Range('E26').Value = '  +1.61%  '

$c = $ws.Range('D27')
$c.Formula = '="166.29"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E27').Value = '  +0.04%  '

$c = $ws.Range('D28')
$c.Formula = '="0.131"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E28').Value = '  +6.32%  '

$c = $ws.Range('D29')
$c.Formula = '="19.04"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E29').Value = '  +1.83%  '

$c = $ws.Range('D30')
$c.Formula = '="1.35"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E30').Value = '  -0.02%  '

$ws.Range('E31').Value = '  +1.81%  '

$c = $ws.Range('D32')
$c.Formula = '="4.54"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E32').Value = '  +0.59%  '

$c = $ws.Range('D33')
$c.Formula = '="4.60"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E33').Value = '  +4.29%  '

$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('E35').Value = '  +8.01%  '

$ws.Range('E36').Value = '  -0.13%  '

$c = $ws.Range('D37')
$c.Formula = '="6.00"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E37').Value = '  +13.54%  '

$c = $ws.Range('D38')
$c.Formula = '="3.31"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E38').Value = '  +4.88%  '

$ws.Range('E39').Value = '  -0.17%  '

$c = $ws.Range('D40')
$c.Formula = '="98.42"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E40').Value = '  +3.77%  '

$ws.Range('E41').Value = '  +1.77%  '

$c = $ws.Range('D42')
$c.Formula = '="1.480.09"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.02%  '

$c = $ws.Range('D43')
$c.Formula = '="0.0940"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E43').Value = '  +2.73%  '

$c = $ws.Range('D45')
$c.Formula = '="16.62"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E45').Value = '  +2.05%  '

$ws.Range('E46').Value = '  +0.37%  '

$ws.Range('E47').Value = '  +16.49%  '

$ws.Range('E48').Value = '  +1.06%  '

$ws.Range('E49').Value = '  +1.54%  '

$c = $ws.Range('D50')
$c.Formula = '="7.07"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E50').Value = '  -0.86%  '

$c = $ws.Range('D51')
$c.Formula = '="2.252.59"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E51').Value = '  +2.32%  '
